$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay text (avoid Excel auto-numeric coercion,
# which would strip trailing zeros / change type vs. the source data).
# Format those cells as Text before writing so COM keeps them as strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.272.61'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.680.68'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.55'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5245'
$ws.Range("E6").Value = '  +2.37%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2702'
$ws.Range("E8").Value = '  +2.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06482'
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.01'
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07529'
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.532'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.664.98'
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5809'
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.68'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.318.52'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.920'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.88'
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.24'
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.50'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.798'
$ws.Range("E25").Value = '  +2.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1244'
$ws.Range("E26").Value = '  +3.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.79'
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06488'
$ws.Range("E28").Value = '  +2.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.357'
$ws.Range("E29").Value = '  +4.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.333'
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.593'
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("E34").Value = '  +1.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6243'
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.404'
$ws.Range("E36").Value = '  +1.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.735'
$ws.Range("E37").Value = '  +3.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.457'
$ws.Range("E38").Value = '  +4.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.111.31'
$ws.Range("E39").Value = '  +2.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01626'
$ws.Range("E40").Value = '  +1.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8769'
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.81'
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.831.81'
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000111'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.04'
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.182'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05274'
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.089'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4292'
$ws.Range("E51").Value = '  -0.03%  '
